# "Estadisticos Segundo Parcial 23 Mayo"
# Fills in the 2nd-partial statistics sheet with real figures, recomputes the
# "Estadisticos Final" sheet to match, and replaces the "Rescatables"
# (students still able to recover subjects) roster with the updated list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Estadisticos 2P" - second partial results (previously a blank template
#    where Blancos==Totales and Reprobados/Aprobados/Por_Apro were all 0)
# ---------------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 6
$ws2P.Range("F2").Value = 30
$ws2P.Range("G2").Value = 83.33
$ws2P.Range("H2").Value = 6.9

$ws2P.Range("D3").Value = 0
$ws2P.Range("E3").Value = 4
$ws2P.Range("F3").Value = 35
$ws2P.Range("G3").Value = 89.73999999999999
$ws2P.Range("H3").Value = 7.3

$ws2P.Range("D4").Value = 0
$ws2P.Range("E4").Value = 5
$ws2P.Range("F4").Value = 19
$ws2P.Range("G4").Value = 79.17
$ws2P.Range("H4").Value = 8

$ws2P.Range("D5").Value = 0
$ws2P.Range("E5").Value = 5
$ws2P.Range("F5").Value = 19
$ws2P.Range("G5").Value = 79.17
$ws2P.Range("H5").Value = 8

$ws2P.Range("D6").Value = 0
$ws2P.Range("E6").Value = 0
$ws2P.Range("F6").Value = 24
$ws2P.Range("G6").Value = 100
$ws2P.Range("H6").Value = 8.1

$ws2P.Range("D7").Value = 0
$ws2P.Range("E7").Value = 0
$ws2P.Range("F7").Value = 24
$ws2P.Range("G7").Value = 100
$ws2P.Range("H7").Value = 8

# ---------------------------------------------------------------------------
# 2) "Estadisticos Final" - updated to the latest (2nd partial) figures;
#    Blancos/Reprobados/Aprobados/Por_Apro now mirror "Estadisticos 2P" while
#    Promedio (H) carries the real running average for the course.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("D2").Value = 0
$wsFinal.Range("E2").Value = 6
$wsFinal.Range("F2").Value = 30
$wsFinal.Range("G2").Value = 83.33
$wsFinal.Range("H2").Value = 7.2

$wsFinal.Range("D3").Value = 0
$wsFinal.Range("E3").Value = 4
$wsFinal.Range("F3").Value = 35
$wsFinal.Range("G3").Value = 89.73999999999999
$wsFinal.Range("H3").Value = 7.4

$wsFinal.Range("D4").Value = 0
$wsFinal.Range("E4").Value = 5
$wsFinal.Range("F4").Value = 19
$wsFinal.Range("G4").Value = 79.17
$wsFinal.Range("H4").Value = 7.6

$wsFinal.Range("D5").Value = 0
$wsFinal.Range("E5").Value = 5
$wsFinal.Range("F5").Value = 19
$wsFinal.Range("G5").Value = 79.17
$wsFinal.Range("H5").Value = 7.6

$wsFinal.Range("D6").Value = 0
$wsFinal.Range("E6").Value = 0
$wsFinal.Range("F6").Value = 24
$wsFinal.Range("G6").Value = 100
$wsFinal.Range("H6").Value = 8.1

$wsFinal.Range("D7").Value = 0
$wsFinal.Range("E7").Value = 0
$wsFinal.Range("F7").Value = 24
$wsFinal.Range("G7").Value = 100
$wsFinal.Range("H7").Value = 8.4

# ---------------------------------------------------------------------------
# 3) "Rescatables" - roster of students still able to recover subjects.
#    The old list (18 students) is replaced by a fresh list of 10 students
#    (rows 12:19 no longer exist).
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

# Drop the now-unused tail rows first so the sheet's dimension becomes A1:G11.
$wsResc.Range("A12:G19").EntireRow.Delete()

$rescatables = @(
    @(2,  24330051920233, "ARELLANO", "PAZ",     "ADRIAN",           "Pensamiento matemático II", "2BEM", 4),
    @(3,  24330051920353, "GARCIA",   "SANCHEZ", "JOY JARA",         "Pensamiento matemático II", "2BEM", 4),
    @(4,  24330051920404, "PARADA",   "SANTOS",  "MARCO DIDIEL",     "Pensamiento matemático II", "2BEM", 4),
    @(5,  23330051920005, "CASTRO",   "ARIAS",   "OMAR DAVID",       "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO", "4AEM", 3),
    @(6,  23330051920005, "CASTRO",   "ARIAS",   "OMAR DAVID",       "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)", "4AEM", 3),
    @(7,  23330051920018, "RAMOS",    "UTRERA",  "CARLOS DAVID",     "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO", "4AEM", 3),
    @(8,  23330051920018, "RAMOS",    "UTRERA",  "CARLOS DAVID",     "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)", "4AEM", 3),
    @(9,  24330051920340, "CASTILLO", "GONZALEZ","ANGEL ALBERTO",    "Pensamiento matemático II", "2BEM", 2),
    @(10, 24330051920324, "FLORES",   "TINOCO",  "ANGEL MOISES",     "Pensamiento matemático II", "2AEM", 1),
    @(11, 24330051920334, "PEREZ",    "ISLAS",   "VICTOR ALEJANDRO", "Pensamiento matemático II", "2BEM", 1)
)

foreach ($row in $rescatables) {
    $r = $row[0]
    $wsResc.Range("A$r").Value = $row[1]
    $wsResc.Range("B$r").Value = $row[2]
    $wsResc.Range("C$r").Value = $row[3]
    $wsResc.Range("D$r").Value = $row[4]
    $wsResc.Range("E$r").Value = $row[5]
    $wsResc.Range("F$r").Value = $row[6]
    $wsResc.Range("G$r").Value = $row[7]
}
